# Add a new calculated column "LEN_LastName" in column J that computes
# the length of each employee's LastName (column C) using an array
# formula that spills from J2 down to J10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("J1").Value = "LEN_LastName"

# Array formula over the whole target range so it is stored as a single
# array formula anchored at J2 with ref J2:J10 (matching LEN(C2:C10)
# being entered once and filled down as an array formula): J2 carries
# <f t="array" ref="J2:J10">LEN(C2:C10)</f> and J3:J10 carry the spilled
# cached values only, exactly like the target workbook.
$ws.Range("J2:J10").FormulaArray = "=LEN(C2:C10)"
